$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1 ("o_10") - update header E1 and refresh row 2 data,
# replacing the old 18-node prompt/solution/response with the new
# 15-node scenario and adding the evaluator_partial_correctness col.
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)

$prompt15 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 15 nodes labelled A to O. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node J? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O
 A 0 4 4 0 0 0 0 0 0 0 0 0 1 0 0
 B 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 4 0 0 1 0 0 2 0 0 0 0 0 0 0 0
 D 0 0 1 0 3 5 0 0 0 0 0 0 0 0 0
 E 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 2 0 0 0 0 4 3 0 0 0 0 3 0
 H 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 3 0 0 5 0 0 0 0 4
 J 0 0 0 0 0 0 0 0 5 0 1 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 2 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 2 0 4 0 0
 M 1 0 0 0 0 0 0 0 0 0 0 4 0 0 0
 N 0 0 0 0 0 0 3 0 0 0 0 0 0 0 2
 O 0 0 0 0 0 0 0 0 4 0 0 0 0 2 0
    
"@

$ws1.Range("A2").Value = $prompt15
$ws1.Range("B2").Value = "A -> M -> L -> K -> J"
$ws1.Range("C2").Value = "The least cost path from node A to node J is A - M - L - K - J."
$ws1.Range("D2").Value = "invalid input"
$ws1.Range("E2").Value = "5/5"
$ws1.Range("A2").EntireRow.AutoFit()

# -----------------------------------------------------------------
# Sheet 2 ("o_20") - new sheet, same header layout, new data row
# for the 25-node scenario.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

$prompt25 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node T? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 1 0 0 5 0
 C 0 1 0 5 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 2 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 4 0 1 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 3 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 3 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 2 0 0 0 0 0 0 3 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 4 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 2 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 5 0 0 0 0 4
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0
 U 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 X 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 4 0 0 0 0 0 0
    
"@

$ws2.Range("A2").Value = $prompt25
$ws2.Range("B2").Value = "A -> B -> O -> S -> T"
$ws2.Range("C2").Value = "The least cost path from node A to node T is A-B-O-S-T with a total cost of 1+2+1+5 = 9."
$ws2.Range("D2").Value = "invalid input"
$ws2.Range("E2").Value = "5/5"
$ws2.Range("A2").EntireRow.AutoFit()

# -----------------------------------------------------------------
# Sheet 3 ("o_20_jumbled") - new sheet, same header layout, new data
# row for the 26-node (jumbled) scenario.
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("D1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)

$prompt26 = @"
 Given is the adjacency matrix for a weighted undirected graph containing 26 nodes labelled A to Z. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   



what is the least cost path from node A to node T? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T U V W X Y Z
 A 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 4 0 3 2 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 1 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 F 0 0 0 0 2 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 4 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 2 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 5 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1
 K 0 0 0 0 0 0 0 0 0 1 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 5 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 3 0 5 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 5 0 4 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 1 0 1 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 1 0 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 U 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 2 0 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 1 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 3 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 4
 Z 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0
    
"@

$ws3.Range("A2").Value = $prompt26
$ws3.Range("B2").Value = "A -> B -> E -> F -> G -> I -> J -> K -> L -> M -> S -> T"
$ws3.Range("C2").Value = "The least cost path from node A to node T is A -> B -> E -> U -> S -> T."
$ws3.Range("D2").Value = "invalid input"
$ws3.Range("E2").Value = "5/12"
$ws3.Range("A2").EntireRow.AutoFit()

$ws1.Activate()
